# :bug: fix: tabela viagens
# The "eventosSazonais" table had its start/end dates for 2024/2025 events;
# this corrects them back to the intended 2023 dates (same month/day, year
# shifted to 2023) and restores the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, start-date serial, end-date serial (OA date serials, matching the
# canonical OOXML <v> values in the corrected workbook)
$rows = @(
    @(2,  44967, 44971),  # Carnaval
    @(3,  44993, 44993),  # Dia Internacional da Mulher
    @(4,  45000, 45000),  # Dia do Consumidor
    @(5,  45033, 45033),  # Páscoa
    @(6,  45054, 45054),  # Dia das Mães
    @(7,  45089, 45089),  # Dia dos Namorados
    @(8,  45101, 45101),  # Dia de São João
    @(9,  45105, 45105),  # Dia Internacional do Orgulho LGBTQIA+
    @(10, 45133, 45133),  # Dia dos Avós
    @(11, 45149, 45149),  # Dia dos Pais
    @(12, 45211, 45211),  # Dia das Crianças
    @(13, 45230, 45230),  # Halloween
    @(14, 45259, 45259),  # Black Friday
    @(15, 45285, 45285),  # Natal
    @(16, 44927, 44927)   # Reveillon
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Restore the selected/active cell left by the editor
$ws.Range("H12").Select()
